$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42647.680914351855
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = $true

$ws.Range("C5").Value = 10162.42
$ws.Range("D5").Value = 10031.01
$ws.Range("E5").Value = 78.63
$ws.Range("F5").Value = 77.599999999999994

$ws.Range("G5").Value = $true
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"

$ws.Range("H5").Value = -1.31

$ws.Range("I5").Value = $true
